# Add "Collections" and "Distributions" report tabs to the Sites report
# workbook, cloning the existing "Sites" tab layout for Collections and
# extending it with a few extra columns for Distributions.

$wb = $excel.ActiveWorkbook

$sites = $wb.Worksheets.Item("Sites")
$individuals = $wb.Worksheets.Item("Individuals")

# --- Sites: deselect it as the active tab but keep the whole-sheet
#     selection behind, while it is still the active sheet. ---
$sites.Range("A1:XFD3").Select()

# --- Collections: clone the Sites tab layout. ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sites.Copy($null, $lastSheet)
$collections = $wb.Worksheets.Item($wb.Worksheets.Count)
$collections.Name = "Collections"

# Materialize the (empty, borderless) column-style placeholder cells under
# the date-range row, matching what Excel leaves behind once that row's
# cells have been touched.
$collections.Range("D2:G2").Value = "x"
$collections.Range("D2:G2").Value = ""

# Resize the Collections columns.
$collections.Columns.Item(1).ColumnWidth = 8.833333
$collections.Columns.Item(2).ColumnWidth = 17.0
$collections.Columns.Item(3).ColumnWidth = 12.833333
$collections.Columns.Item(4).ColumnWidth = 12.5
$collections.Columns.Item(5).ColumnWidth = 11.5
$collections.Columns.Item(6).ColumnWidth = 13.666667
$collections.Columns.Item(7).ColumnWidth = 12.666667
$collections.Columns.Item(8).ColumnWidth = 10.5
$collections.Columns.Item(9).ColumnWidth = 5.5
$collections.Columns.Item(10).ColumnWidth = 8.833333
$collections.Columns.Item(11).ColumnWidth = 10.333333

$collections.Range("B10").Select()

# --- Distributions: clone the Collections tab layout, then append the
#     extra release / acclimation columns. ---
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$collections.Copy($null, $lastSheet2)
$distributions = $wb.Worksheets.Item($wb.Worksheets.Count)
$distributions.Name = "Distributions"

$distributions.Range("L3").Value = "Release Method"
$distributions.Range("M3").Value = "Lifestage"
$distributions.Range("N3").Value = "Truck Temp"
$distributions.Range("O3").Value = "River Temp"
$distributions.Range("P3").Value = "Acclimation Time (mins)"

# L3/M3 use the same header look as the rest of the row; N3:P3 borrow the
# other accent header style already used on the Individuals tab.
$individuals.Range("F3").Copy()
$distributions.Range("L3:M3").PasteSpecial(-4122)
$individuals.Range("A3").Copy()
$distributions.Range("N3:P3").PasteSpecial(-4122)

$distributions.Columns.Item(12).ColumnWidth = 14.666667
$distributions.Columns.Item(13).ColumnWidth = 8.166667
$distributions.Columns.Item(14).ColumnWidth = 13.166667
$distributions.Columns.Item(15).ColumnWidth = 12.333333
$distributions.Columns.Item(16).ColumnWidth = 22.0

# Also re-tune the columns Distributions inherited from Collections, since
# the Distributions tab needs slightly narrower C/E/G/J/K columns.
$distributions.Columns.Item(3).ColumnWidth = 4.333333
$distributions.Columns.Item(5).ColumnWidth = 9.166667
$distributions.Columns.Item(7).ColumnWidth = 9.666667
$distributions.Columns.Item(10).ColumnWidth = 6.166667
$distributions.Columns.Item(11).ColumnWidth = 9.0

$distributions.Range("L3").Select()

$collections.Activate()
